$d = $word.ActiveDocument

# Commit: "Use Assessment of Significance and tweak carriage return"
# The paragraph holding the "<Assessment of Significance>" placeholder is
# immediately followed by a blank paragraph (just a carriage return). That
# extra blank paragraph is removed, merging its paragraph mark away.

$range = $d.Content
$range.Find.ClearFormatting()
$found = $range.Find.Execute(
    "<Assessment of Significance>", $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if ($found) {
    # $range now spans the found text; the paragraph that contains it is the
    # anchor paragraph, and the next paragraph is the blank one to remove.
    $anchorPara = $range.Paragraphs(1)
    $blankPara = $anchorPara.Next()

    if ($blankPara.Range.Text -eq "`r") {
        $blankPara.Range.Delete()
        Write-Output "Removed blank paragraph after 'Assessment of Significance'."
    } else {
        Write-Output "Unexpected content after 'Assessment of Significance'; no change made."
    }
} else {
    Write-Output "Could not find '<Assessment of Significance>'."
}
